$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: refresh crypto price/volume data.
# NumberFormat '@' forces text interpretation so numeric-looking strings
# (e.g. '0.996') are not silently converted to numbers; ClearFormats()
# afterwards removes the style footprint again so no stray style id is left.
$cellValues = @{
    'D2' = '56.949.62'
    'E2' = '  +2.56%  '
    'D3' = '2.507.85'
    'E3' = '  +0.25%  '
    'E4' = '  +0.08%  '
    'D5' = '497.08'
    'E5' = '  +3.10%  '
    'D6' = '154.26'
    'E6' = '  +9.92%  '
    'D7' = '0.996'
    'E7' = '  -0.31%  '
    'D8' = '0.516'
    'E8' = '  +1.39%  '
    'D9' = '2.521.23'
    'E9' = '  +0.86%  '
    'D10' = '5.79'
    'E10' = '  +5.78%  '
    'D11' = '0.0995'
    'E11' = '  +0.90%  '
    'E12' = '  +3.16%  '
    'E13' = '  +1.39%  '
    'D14' = '2.931.88'
    'E14' = '  -0.04%  '
    'D15' = '57.047.04'
    'E15' = '  +2.72%  '
    'D16' = '21.58'
    'E16' = '  +4.46%  '
    'D17' = '0.0000138'
    'E17' = '  +0.05%  '
    'D18' = '2.521.14'
    'E18' = '  +0.99%  '
    'D19' = '4.57'
    'E19' = '  +4.52%  '
    'D20' = '10.39'
    'E20' = '  +3.51%  '
    'D21' = '326.15'
    'E21' = '  +1.82%  '
    'D22' = '0.998'
    'E22' = '  -0.03%  '
    'E23' = '  +4.27%  '
    'D24' = '59.14'
    'E24' = '  +2.29%  '
    'E25' = '  +0.94%  '
    'E26' = '  -2.08%  '
    'D27' = '0.999'
    'E27' = '  -1.19%  '
    'D28' = '2.612.57'
    'E28' = '  +0.30%  '
    'D29' = '7.70'
    'E29' = '  +4.49%  '
    'D30' = '0.0₃0821'
    'E30' = '  +2.99%  '
    'D31' = '0.998'
    'D32' = '151.54'
    'E32' = '  +1.43%  '
    'D33' = '18.46'
    'E33' = '  +1.61%  '
    'E34' = '  +3.83%  '
    'E35' = '  +2.20%  '
    'D37' = '3.81'
    'E37' = '  +2.62%  '
    'D38' = '0.880'
    'E38' = '  +2.17%  '
    'D39' = '1.40'
    'E39' = '  +6.46%  '
    'D40' = '34.42'
    'E40' = '  +0.64%  '
    'D41' = '0.0569'
    'E41' = '  +2.64%  '
    'E42' = '  +3.67%  '
    'D43' = '0.617'
    'E43' = '  +1.42%  '
    'D44' = '0.994'
    'E44' = '  -0.56%  '
    'D45' = '4.99'
    'E45' = '  +7.97%  '
    'D46' = '270.49'
    'E46' = '  +7.62%  '
    'D47' = '0.0933'
    'E47' = '  +3.28%  '
    'D48' = '0.0232'
    'E48' = '  +3.94%  '
    'D49' = '10.22'
    'E49' = '  +0.65%  '
    'D50' = '17.96'
    'E50' = '  +2.06%  '
    'D51' = '1.915.65'
    'E51' = '  -3.25%  '
}

foreach ($ref in $cellValues.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $cellValues[$ref]
    $cell.ClearFormats()
}
